$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''60.814.01'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '''  -1.23%  '
$ws.Range('E2').ClearFormats()
$ws.Range('D3').Value = '''3.376.41'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '''  -0.48%  '
$ws.Range('E3').ClearFormats()
$ws.Range('E4').Value = '''  -0.14%  '
$ws.Range('E4').ClearFormats()
$ws.Range('D5').Value = '''570.10'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '''  -1.25%  '
$ws.Range('E5').ClearFormats()
$ws.Range('D6').Value = '''135.93'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '''  -0.81%  '
$ws.Range('E6').ClearFormats()
$ws.Range('E7').Value = '''  +0.00%  '
$ws.Range('E7').ClearFormats()
$ws.Range('D8').Value = '''3.372.37'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '''  -0.59%  '
$ws.Range('E8').ClearFormats()
$ws.Range('D9').Value = '''0.468'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '''  -1.12%  '
$ws.Range('E9').ClearFormats()
$ws.Range('D10').Value = '''7.57'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '''  +0.91%  '
$ws.Range('E10').ClearFormats()
$ws.Range('E11').Value = '''  -3.25%  '
$ws.Range('E11').ClearFormats()
$ws.Range('E12').Value = '''  -2.98%  '
$ws.Range('E12').ClearFormats()
$ws.Range('D13').Value = '''3.951.35'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '''  -0.62%  '
$ws.Range('E13').ClearFormats()
$ws.Range('E14').Value = '''  -0.62%  '
$ws.Range('E14').ClearFormats()
$ws.Range('D15').Value = '''25.96'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '''  +0.22%  '
$ws.Range('E15').ClearFormats()
$ws.Range('D16').Value = '''3.375.98'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '''  -0.57%  '
$ws.Range('E16').ClearFormats()
$ws.Range('D17').Value = '''0.0000170'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '''  -4.05%  '
$ws.Range('E17').ClearFormats()
$ws.Range('D18').Value = '''60.869.54'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '''  -1.33%  '
$ws.Range('E18').ClearFormats()
$ws.Range('D19').Value = '''5.82'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '''  -1.25%  '
$ws.Range('E19').ClearFormats()
$ws.Range('D20').Value = '''13.73'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '''  -3.24%  '
$ws.Range('E20').ClearFormats()
$ws.Range('D21').Value = '''9.23'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '''  -2.17%  '
$ws.Range('E21').ClearFormats()
$ws.Range('D22').Value = '''371.94'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '''  -1.27%  '
$ws.Range('E22').ClearFormats()
$ws.Range('D23').Value = '''3.509.09'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '''  -0.71%  '
$ws.Range('E23').ClearFormats()
$ws.Range('D24').Value = '''0.547'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '''  -2.02%  '
$ws.Range('E24').ClearFormats()
$ws.Range('E25').Value = '''  +0.12%  '
$ws.Range('E25').ClearFormats()
$ws.Range('D26').Value = '''70.72'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '''  -0.65%  '
$ws.Range('E26').ClearFormats()
$ws.Range('D27').Value = '''0.0000123'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '''  -3.09%  '
$ws.Range('E27').ClearFormats()
$ws.Range('D28').Value = '''0.175'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '''  +8.83%  '
$ws.Range('E28').ClearFormats()
$ws.Range('E29').Value = '''  -5.42%  '
$ws.Range('E29').ClearFormats()
$ws.Range('D30').Value = '''0.999'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '''  -0.05%  '
$ws.Range('E30').ClearFormats()
$ws.Range('D31').Value = '''7.33'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '''  -2.66%  '
$ws.Range('E31').ClearFormats()
$ws.Range('E32').Value = '''  -3.06%  '
$ws.Range('E32').ClearFormats()
$ws.Range('E33').Value = '''  -2.67%  '
$ws.Range('E33').ClearFormats()
$ws.Range('E35').Value = '''  -0.95%  '
$ws.Range('E35').ClearFormats()
$ws.Range('E36').Value = '''  -4.45%  '
$ws.Range('E36').ClearFormats()
$ws.Range('E37').Value = '''  -1.39%  '
$ws.Range('E37').ClearFormats()
$ws.Range('D38').Value = '''6.76'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '''  -1.08%  '
$ws.Range('E38').ClearFormats()
$ws.Range('D39').Value = '''164.77'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '''  -0.62%  '
$ws.Range('E39').ClearFormats()
$ws.Range('D40').Value = '''0.0757'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '''  -2.96%  '
$ws.Range('E40').ClearFormats()
$ws.Range('D41').Value = '''1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('D42').Value = '''1.72'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '''  -0.77%  '
$ws.Range('E42').ClearFormats()
$ws.Range('E43').Value = '''  -1.37%  '
$ws.Range('E43').ClearFormats()
$ws.Range('B44').Value = '''OKB'
$ws.Range('B44').ClearFormats()
$ws.Range('C44').Value = '''https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('C44').ClearFormats()
$ws.Range('D44').Value = '''41.86'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '''  +0.94%  '
$ws.Range('E44').ClearFormats()
$ws.Range('B45').Value = '''EnergySwap'
$ws.Range('B45').ClearFormats()
$ws.Range('C45').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C45').ClearFormats()
$ws.Range('D45').Value = '''25.07'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '''  -0.63%  '
$ws.Range('E45').ClearFormats()
$ws.Range('D46').Value = '''4.32'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '''  -2.05%  '
$ws.Range('E46').ClearFormats()
$ws.Range('E47').Value = '''  -6.54%  '
$ws.Range('E47').ClearFormats()
$ws.Range('D48').Value = '''2.529.28'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '''  +8.20%  '
$ws.Range('E48').ClearFormats()
$ws.Range('D49').Value = '''23.47'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '''  +3.38%  '
$ws.Range('E49').ClearFormats()
$ws.Range('D50').Value = '''6.75'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '''  -1.55%  '
$ws.Range('E50').ClearFormats()
$ws.Range('E51').Value = '''  +1.07%  '
$ws.Range('E51').ClearFormats()
